$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width changes: B and C (81 -> 70) ---
$ws.Range("B1").ColumnWidth = 69.17
$ws.Range("C1").ColumnWidth = 69.17

# --- Column width changes for U..AL (21..38) ---
$widths = @{
    "U"  = 34.17;  "V"  = 36.17;  "W"  = 41.17;  "X"  = 43.17
    "Y"  = 31.17;  "Z"  = 33.17;  "AA" = 28.17;  "AB" = 30.17
    "AC" = 37.17;  "AD" = 39.17;  "AE" = 23.17;  "AF" = 25.17
    "AG" = 24.17;  "AH" = 26.17;  "AI" = 27.17;  "AJ" = 29.17
    "AK" = 24.17;  "AL" = 26.17
}
foreach ($col in $widths.Keys) {
    $ws.Range($col + "1").ColumnWidth = $widths[$col]
}

# --- Header row (row 1): rename/re-order the link_executionDetails / link_executionLinks
#     groups into link_executionLinks (moved earlier) + new link_jobDetails group ---
$headers = @{
    "U1"  = "link_executionLinks_executions_id"
    "V1"  = "link_executionLinks_executions_id_1"
    "W1"  = "link_executionLinks_internalRoleLinkName"
    "X1"  = "link_executionLinks_internalRoleLinkName_1"
    "Y1"  = "link_executionLinks_project_id"
    "Z1"  = "link_executionLinks_project_id_1"
    "AA1" = "link_executionLinks_team_id"
    "AB1" = "link_executionLinks_team_id_1"
    "AC1" = "link_jobDetails_internalRoleLinkName"
    "AD1" = "link_jobDetails_internalRoleLinkName_1"
    "AE1" = "link_jobDetails_job_id"
    "AF1" = "link_jobDetails_job_id_1"
    "AG1" = "link_jobDetails_plan_id"
    "AH1" = "link_jobDetails_plan_id_1"
    "AI1" = "link_jobDetails_project_id"
    "AJ1" = "link_jobDetails_project_id_1"
    "AK1" = "link_jobDetails_team_id"
    "AL1" = "link_jobDetails_team_id_1"
}
foreach ($addr in $headers.Keys) {
    $ws.Range($addr).Value = $headers[$addr]
}

# --- Data row (row 2): values for the same columns, following the new header order ---
# The values are numeric-looking text (not real numbers) in the source data, so we stage
# them in a scratch range formatted as Text, then paste-special (values only) into the
# target cells -- this keeps the cells' number storage as text without leaving the
# Text number-format behind on the destination cells.
$scratch = "A200:R200"
$ws.Range($scratch).NumberFormat = "@"
$ws.Range("A200").Value = "10"
$ws.Range("B200").Value = "12"
$ws.Range("C200").Value = "10"
$ws.Range("D200").Value = "12"
$ws.Range("E200").Value = "1588984"
$ws.Range("F200").Value = "1588984"
$ws.Range("G200").Value = "1570311"
$ws.Range("H200").Value = "1570311"
$ws.Range("I200").Value = "8"
$ws.Range("J200").Value = "10"
$ws.Range("K200").Value = "8"
$ws.Range("L200").Value = "10"
$ws.Range("M200").Value = "837097"
$ws.Range("N200").Value = "837132"
$ws.Range("O200").Value = "1588984"
$ws.Range("P200").Value = "1588984"
$ws.Range("Q200").Value = "1570311"
$ws.Range("R200").Value = "1570311"

$ws.Range($scratch).Copy()
$ws.Range("U2").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()

# --- B2 / C2 text update ---
$ws.Range("B2").Value = "Data Files/AI-Generated/Common/scheduleAndRunTestExecution-test-data"
$ws.Range("C2").Value = "Data Files/AI-Generated/Common/scheduleAndRunTestExecution-test-data"
